$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mysql")
$ws.Range("B2").Value = "select count(*) from orderadataplannew"
